$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.296.68'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '2.656.14'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("D5").Value = "'597.58"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").Value = "'175.00"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").Value = '2.656.30'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("D13").Value = "'5.00"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = '3.140.89'
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '72.105.02'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").Value = "'26.26"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = '2.644.60'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = "'12.23"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +5.42%  '
$ws.Range("D20").Value = "'8.15"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("D21").Value = "'370.63"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -2.51%  '
$ws.Range("D22").Value = "'4.19"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = "'2.05"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").Value = "'72.14"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").Value = "'4.31"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -1.78%  '
$ws.Range("D27").Value = "'9.75"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -2.34%  '
$ws.Range("D28").Value = '2.792.31'
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = '0.0₃0970'
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = "'498.23"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -4.05%  '
$ws.Range("D33").Value = "'1.30"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").Value = "'19.51"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("E38").Value = '  -0.63%  '
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("E40").Value = '  -1.60%  '
$ws.Range("D41").Value = "'1.77"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -4.36%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = "'5.00"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").Value = "'155.42"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +3.86%  '
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("D49").Value = "'0.557"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("D50").Value = "'1.73"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +1.97%  '
$ws.Range("D51").Value = "'0.0757"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -1.14%  '
